$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New content cells (order matters for shared-string index allocation) ---
$ws.Range("I28").Value = "SELECT"
$ws.Range("I29").Value = "FROM"
$ws.Range("I30").Value = "WHERE"
$ws.Range("I31").Value = "GROUP BY"
$ws.Range("I32").Value = "HAVING"
$ws.Range("I33").Value = "ORDER BY"
$ws.Range("I34").Value = "LIMIT"
$ws.Range("L15").Value = "runner"
$ws.Range("M16").Value = "avg()"
$ws.Range("I27").Value = "What is the successful delivery percentage for each runner?"
$ws.Range("M15").Value = "successful_perc"

# --- Highlight the pizza_id = 1 (Meatlovers) rows in customer_orders with yellow fill ---
$ws.Range("D11").Interior.Color = 65535
$ws.Range("D12").Interior.Color = 65535
$ws.Range("D13").Interior.Color = 65535
$ws.Range("D15").Interior.Color = 65535
$ws.Range("D16").Interior.Color = 65535
$ws.Range("D18").Interior.Color = 65535
$ws.Range("D21").Interior.Color = 65535
$ws.Range("D22").Interior.Color = 65535
$ws.Range("D23").Interior.Color = 65535
$ws.Range("D24").Interior.Color = 65535

# --- Column widths for the newly used columns L, M, N ---
$ws.Columns.Item(12).ColumnWidth = 16.5
$ws.Columns.Item(13).ColumnWidth = 16.333333333333332
$ws.Columns.Item(14).ColumnWidth = 21

# --- Update the active selection / view ---
$ws.Range("L12").Select()
